# #CRM-1615 Add Remarks in SF pending booking page
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Remarks" column header (Q1) + merge-field placeholder (Q2)
$ws.Range("P1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q1").Value = "Remarks"

$ws.Range("H2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q2").Value = "{bookings:booking_remarks}"

$excel.CutCopyMode = $false

# Restore the view: scrolled so column G is leftmost, P4 selected
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("P4").Select() | Out-Null
